$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "CHANGE" + "S:" (two runs) -> single run "CHANGES:"
#    Find/Replace across the run boundary collapses the matched text into a
#    single run that keeps the first run's formatting (matches target XML).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("CHANGES:", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "CHANGES:", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Insert a new bulleted line right after the
#    "[Light support version] -> maximum number of blocks lowered to 20000."
#    paragraph, reusing that paragraph's formatting (same NoSpacing style +
#    numPr list) by duplicating the whole paragraph (incl. its mark) and then
#    swapping in the new run content (with proofErr wrappers) via WordOpenXML.
# ---------------------------------------------------------------------------
$maxBlocksPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*maximum number of blocks lowered to 20000.*") {
        $maxBlocksPara = $p
    }
}

$maxBlocksPara.Range.Copy()
$afterMaxBlocks = $d.Range($maxBlocksPara.Range.End, $maxBlocksPara.Range.End)
$afterMaxBlocks.Paste()

# The pasted duplicate is the paragraph right after the original one.
$newPara = $maxBlocksPara.Next()
$newInner = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)

$rPr = '<w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr>'
$runs = ""
$runs += "<w:r>$rPr<w:t>[Light support version]</w:t></w:r>"
$runs += "<w:r>$rPr<w:t xml:space=`"preserve`"> -&gt; </w:t></w:r>"
$runs += '<w:proofErr w:type="gramStart"/>'
$runs += "<w:r>$rPr<w:t>water</w:t></w:r>"
$runs += '<w:proofErr w:type="gramEnd"/>'
$runs += "<w:r>$rPr<w:t xml:space=`"preserve`"> effects are </w:t></w:r>"
$runs += '<w:proofErr w:type="spellStart"/>'
$runs += "<w:r>$rPr<w:t>not longer</w:t></w:r>"
$runs += '<w:proofErr w:type="spellEnd"/>'
$runs += "<w:r>$rPr<w:t xml:space=`"preserve`"> supported.</w:t></w:r>"

$newParaXml = '<?xml version="1.0" standalone="yes"?>' + `
    '<?mso-application progid="Word.Document"?>' + `
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    "<w:body><w:p>$runs</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"

$newInner.InsertXML($newParaXml)

# ---------------------------------------------------------------------------
# 3) Strip the explicit color/size formatting from the paragraph-mark rPr of
#    the blank "NoSpacing" paragraph right after "...now fixed." (keep only
#    rFonts).
# ---------------------------------------------------------------------------
$chunkPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Chunk function which check in which chunk is certain block is now fixed.*") {
        $chunkPara = $p
    }
}
$blankAfterChunk = $chunkPara.Next()
$blankFull = $d.Range($blankAfterChunk.Range.Start, $blankAfterChunk.Range.End)
$blankXml = '<?xml version="1.0" standalone="yes"?>' + `
    '<?mso-application progid="Word.Document"?>' + `
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:body><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$blankFull.InsertXML($blankXml)

# ---------------------------------------------------------------------------
# 4) Remove the empty "ind left=360" spacer paragraph (merging it away), then
#    merge the bookmark-only paragraph with the following "Available
#    commands:" paragraph (taking the latter's bold/48pt formatting), and
#    re-create the _GoBack bookmark that anchors the merge point.
# ---------------------------------------------------------------------------
$availPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Available commands:*") {
        $availPara = $p
    }
}
$bookmarkPara = $availPara.Previous()
$spacerPara = $bookmarkPara.Previous()

# Delete the spacer paragraph entirely (text + its own paragraph mark).
$spacerPara.Range.Delete()

# Delete the bookmark paragraph's own mark so it merges into "Available
# commands:", which then keeps its own (bold, 48pt) paragraph formatting.
$bookmarkMark = $d.Range($bookmarkPara.Range.End - 1, $bookmarkPara.Range.End)
$bookmarkMark.Delete()

# Re-insert the _GoBack bookmark at the (now collapsed) original location.
$mergedPara = $availPara
$bmRange = $d.Range($mergedPara.Range.Start, $mergedPara.Range.Start)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
